$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'55.595.00"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +7.83%  "
$ws.Range("D3").Value = "'2.502.70"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +9.99%  "
$ws.Range("D5").Value = "'486.20"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +13.84%  "
$ws.Range("D6").Value = "'141.00"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +19.97%  "
$ws.Range("E7").Value = "  +0.28%  "
$ws.Range("D8").Value = "'0.512"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +12.25%  "
$ws.Range("D9").Value = "'2.501.58"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +11.25%  "
$ws.Range("D10").Value = "'0.0988"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +15.36%  "
$ws.Range("D11").Value = "'5.53"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +8.61%  "
$ws.Range("D12").Value = "'0.330"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +12.43%  "
$ws.Range("E13").Value = "  +2.00%  "
$ws.Range("D14").Value = "'2.944.35"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +12.83%  "
$ws.Range("D15").Value = "'55.667.67"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +7.83%  "
$ws.Range("D16").Value = "'20.78"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +13.42%  "
$ws.Range("D17").Value = "'0.0000138"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +20.16%  "
$ws.Range("D18").Value = "'2.512.86"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +10.62%  "
$ws.Range("D19").Value = "'4.39"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +14.51%  "
$ws.Range("D20").Value = "'323.04"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +11.07%  "
$ws.Range("D21").Value = "'10.01"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +15.77%  "
$ws.Range("D22").Value = "'1.00"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.16%  "
$ws.Range("D23").Value = "'5.76"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +12.65%  "
$ws.Range("D24").Value = "'58.32"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +11.25%  "
$ws.Range("D25").Value = "'0.170"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +21.21%  "
$ws.Range("E26").Value = "  +16.76%  "
$ws.Range("E27").Value = "  +0.07%  "
$ws.Range("D28").Value = "'2.627.74"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +10.71%  "
$ws.Range("D29").Value = "'7.46"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +11.37%  "
$ws.Range("D30").Value = "'0.0₃0797"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +21.03%  "
$ws.Range("E31").Value = "  +0.36%  "
$ws.Range("D32").Value = "'150.35"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +7.01%  "
$ws.Range("D33").Value = "'18.29"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +10.02%  "
$ws.Range("D34").Value = "'1.50"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +15.80%  "
$ws.Range("D35").Value = "'5.23"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +14.43%  "
$ws.Range("D36").Value = "'0.874"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +10.59%  "
$ws.Range("D37").Value = "'3.71"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +10.43%  "
$ws.Range("D38").Value = "'1.12"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +15.82%  "
$ws.Range("D39").Value = "'34.30"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +9.49%  "
$ws.Range("D40").Value = "'0.616"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +21.00%  "
$ws.Range("D41").Value = "'0.0558"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +15.15%  "
$ws.Range("E42").Value = "  +0.39%  "
$ws.Range("D43").Value = "'3.44"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +12.67%  "
$ws.Range("D44").Value = "'1.32"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +12.77%  "
$ws.Range("B45").Value = "Maker"
$ws.Range("C45").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D45").Value = "'2.003.18"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +9.04%  "
$ws.Range("B46").Value = "RenderToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D46").Value = "'4.70"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +26.26%  "
$ws.Range("B47").Value = "Stellar"
$ws.Range("C47").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D47").Value = "'0.0914"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +13.52%  "
$ws.Range("B48").Value = "WhiteBITCoin"
$ws.Range("C48").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D48").Value = "'10.12"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.38%  "
$ws.Range("B49").Value = "Bittensor"
$ws.Range("C49").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D49").Value = "'255.26"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +42.46%  "
$ws.Range("B50").Value = "VeChain"
$ws.Range("C50").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D50").Value = "'0.0226"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +12.93%  "
$ws.Range("D51").Value = "'17.56"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +14.92%  "
